$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the static predicted value in C2 with a formula that derives
# "TimeTaken in Hours" from "TimeTaken in Minutes" (B2/60).
$ws.Range("C2").Formula = "=B2/60"
